$d = $word.ActiveDocument

# Right single quotation mark used throughout the document for possessives.
$rq = [char]0x2019

# ---------------------------------------------------------------------------
# 1) Split "... the optimal bid and the bot's previous bid, and a human
#    utility higher ..." so that "(excluding) " is inserted right before
#    "bot's previous bid".
# ---------------------------------------------------------------------------
$old1 = " the optimal bid and the bot" + $rq + "s previous bid, and a human utility higher"
$new1 = " the optimal bid and the (excluding) bot" + $rq + "s previous bid, and a human utility higher"
$found1 = $d.Content.Find.Execute($old1, $true, $false, $false, $false, $false, $true, 1, $false, $new1, 2)

# ---------------------------------------------------------------------------
# 2) Append the new explanatory sentence right after "... the bot's previous
#    bid." (the second occurrence, the one that ends the bullet).
# ---------------------------------------------------------------------------
$old2 = " or equal to the bot" + $rq + "s previous bid."
$new2 = " or equal to the bot" + $rq + "s previous bid. We exclude the bot" + $rq + "s previous bid utility unless that gives us an empty set of possible concessions."
$found2 = $d.Content.Find.Execute($old2, $true, $false, $false, $false, $false, $true, 1, $false, $new2, 2)

# ---------------------------------------------------------------------------
# 3) Move the "_GoBack" bookmark from the end of the document (last, empty
#    paragraph) into the newly inserted sentence, between "empty se" and
#    "t of possible concessions."
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

$old3 = "unless that gives us an empty se"
$rng3 = $d.Content
$found3 = $rng3.Find.Execute($old3, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found3) {
    $gb = $d.Range($rng3.End, $rng3.End)
    $d.Bookmarks.Add("_GoBack", $gb)
}
